# Updated cryptos list on Thu Mar  9 23:35:58 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "20.390.44"
$ws.Range("E2").Value = "  -6.54%  "
$ws.Range("D3").Value = "1.438.32"
$ws.Range("E3").Value = "  -6.71%  "
$ws.Range("E4").Value = "  -0.32%  "
$ws.Range("E5").Value = "  -0.21%  "
$ws.Range("D6").Value = "'277.35"
$ws.Range("E6").Value = "  -4.21%  "
$ws.Range("D7").Value = "'0.3726"
$ws.Range("E7").Value = "  -5.36%  "
$ws.Range("D8").Value = "'0.3086"
$ws.Range("E8").Value = "  -3.60%  "
$ws.Range("E9").Value = "  -6.47%  "
$ws.Range("E10").Value = "  -4.99%  "
$ws.Range("D11").Value = "'0.06583"
$ws.Range("E11").Value = "  -8.30%  "
$ws.Range("D12").Value = "'1.002"
$ws.Range("E12").Value = "  -0.35%  "
$ws.Range("D13").Value = "'5.371"
$ws.Range("E13").Value = "  -4.64%  "
$ws.Range("D14").Value = "'17.32"
$ws.Range("E14").Value = "  -6.82%  "
$ws.Range("D15").Value = "'6.135"
$ws.Range("E15").Value = "  -7.63%  "
$ws.Range("D16").Value = "1.443.27"
$ws.Range("E16").Value = "  -6.65%  "
$ws.Range("E17").Value = "  -7.96%  "
$ws.Range("D18").Value = "'76.56"
$ws.Range("E18").Value = "  -8.15%  "
$ws.Range("D19").Value = "'0.05822"
$ws.Range("E19").Value = "  -11.32%  "
$ws.Range("E20").Value = "  -0.24%  "
$ws.Range("D21").Value = "'5.732"
$ws.Range("E21").Value = "  -6.83%  "
$ws.Range("D22").Value = "'14.41"
$ws.Range("E22").Value = "  -5.68%  "
$ws.Range("D23").Value = "'10.90"
$ws.Range("E23").Value = "  -1.15%  "
$ws.Range("D24").Value = "'2.320"
$ws.Range("E24").Value = "  -2.75%  "
$ws.Range("D25").Value = "20.389.92"
$ws.Range("E25").Value = "  -6.62%  "
$ws.Range("E26").Value = "  -6.19%  "
$ws.Range("D27").Value = "'142.63"
$ws.Range("E27").Value = "  -1.74%  "
$ws.Range("D28").Value = "'17.02"
$ws.Range("E28").Value = "  -7.46%  "
$ws.Range("D29").Value = "1.601.55"
$ws.Range("E29").Value = "  -6.80%  "
$ws.Range("D30").Value = "'110.18"
$ws.Range("E30").Value = "  -5.99%  "
$ws.Range("D31").Value = "'3.975"
$ws.Range("E31").Value = "  -18.01%  "
$ws.Range("D32").Value = "'0.9138"
$ws.Range("E32").Value = "  -5.73%  "
$ws.Range("D33").Value = "'5.482"
$ws.Range("E33").Value = "  -6.82%  "
$ws.Range("D34").Value = "'0.07706"
$ws.Range("E34").Value = "  -6.49%  "
$ws.Range("D35").Value = "'8.363"
$ws.Range("E35").Value = "  -6.76%  "
$ws.Range("D36").Value = "'10.96"
$ws.Range("E36").Value = "  +3.05%  "
$ws.Range("B37").Value = "Hedera"
$ws.Range("C37").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D37").Value = "'0.05723"
$ws.Range("E37").Value = "  -5.85%  "
$ws.Range("B38").Value = "Frax"
$ws.Range("C38").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D38").Value = "'1.002"
$ws.Range("E38").Value = "  -0.21%  "
$ws.Range("B39").Value = "InternetComputer(DFINITY)"
$ws.Range("C39").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D39").Value = "'4.732"
$ws.Range("E39").Value = "  -7.47%  "
$ws.Range("B40").Value = "TrustWalletToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D40").Value = "'1.136"
$ws.Range("E40").Value = "  -5.16%  "
$ws.Range("E41").Value = "  -5.67%  "
$ws.Range("D42").Value = "'0.02028"
$ws.Range("E42").Value = "  -9.27%  "
$ws.Range("D43").Value = "'1.343"
$ws.Range("E43").Value = "  -13.01%  "
$ws.Range("E44").Value = "  -4.39%  "
$ws.Range("D45").Value = "'0.5327"
$ws.Range("E45").Value = "  -7.54%  "
$ws.Range("D46").Value = "'12.10"
$ws.Range("E46").Value = "  -6.55%  "
$ws.Range("D47").Value = "'0.5166"
$ws.Range("E47").Value = "  -6.90%  "
$ws.Range("D48").Value = "'112.31"
$ws.Range("E48").Value = "  -4.02%  "
$ws.Range("D49").Value = "'1.789"
$ws.Range("E49").Value = "  -3.84%  "
$ws.Range("E50").Value = "  -6.50%  "
$ws.Range("D51").Value = "'1.002"
$ws.Range("E51").Value = "  -0.30%  "
